# Correct the IFRS "company_list" figures (삼성카드) for rows 2-9 (data rows 1-8).
# The previous export pulled the wrong quarter/annual columns out of the source
# table, so every financial metric cell for these rows is replaced with the right
# figure, and the handful of columns (J/O/U, and Q:U for rows 7-9) that no longer
# carry a value for this data set are cleared out entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 35218
$ws.Range("E2").Value = 8654
$ws.Range("F2").Value = 8654
$ws.Range("G2").Value = 8644
$ws.Range("H2").Value = 6560
$ws.Range("I2").Value = 6560
$ws.Range("K2").Value = 177366
$ws.Range("L2").Value = 113463
$ws.Range("M2").Value = 63903
$ws.Range("N2").Value = 63903
$ws.Range("P2").Value = 6148
$ws.Range("Q2").Value = -11927
$ws.Range("R2").Value = 4690
$ws.Range("S2").Value = 9084
$ws.Range("T2").Value = 502
$ws.Range("V2").Value = 90957
$ws.Range("W2").Value = 24.57
$ws.Range("X2").Value = 18.63
$ws.Range("Y2").Value = 10.42
$ws.Range("Z2").Value = 3.82
$ws.Range("AA2").Value = 177.55
$ws.Range("AB2").Value = 940.4299999999999
$ws.Range("AC2").Value = 5662
$ws.Range("AD2").Value = 7.81
$ws.Range("AE2").Value = 55375
$ws.Range("AF2").Value = 0.8
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 2.26
$ws.Range("AI2").Value = 17.59
$ws.Range("AJ2").Value = 115858891
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("U2").ClearContents()

# Row 3
$ws.Range("D3").Value = 33022
$ws.Range("E3").Value = 3842
$ws.Range("F3").Value = 3842
$ws.Range("G3").Value = 4257
$ws.Range("H3").Value = 3337
$ws.Range("I3").Value = 3337
$ws.Range("K3").Value = 190710
$ws.Range("L3").Value = 123827
$ws.Range("M3").Value = 66883
$ws.Range("N3").Value = 66883
$ws.Range("P3").Value = 6148
$ws.Range("Q3").Value = -4041
$ws.Range("R3").Value = -866
$ws.Range("S3").Value = 5218
$ws.Range("T3").Value = 319
$ws.Range("V3").Value = 98254
$ws.Range("W3").Value = 11.63
$ws.Range("X3").Value = 10.11
$ws.Range("Y3").Value = 5.1
$ws.Range("Z3").Value = 1.81
$ws.Range("AA3").Value = 185.14
$ws.Range("AB3").Value = 988.9
$ws.Range("AC3").Value = 2880
$ws.Range("AD3").Value = 10.71
$ws.Range("AE3").Value = 57957
$ws.Range("AF3").Value = 0.53
$ws.Range("AG3").Value = 1500
$ws.Range("AH3").Value = 4.86
$ws.Range("AI3").Value = 51.87
$ws.Range("AJ3").Value = 115858891
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("U3").ClearContents()

# Row 4
$ws.Range("D4").Value = 34701
$ws.Range("E4").Value = 4309
$ws.Range("F4").Value = 4309
$ws.Range("G4").Value = 4563
$ws.Range("H4").Value = 3494
$ws.Range("I4").Value = 3494
$ws.Range("K4").Value = 219045
$ws.Range("L4").Value = 152888
$ws.Range("M4").Value = 66157
$ws.Range("N4").Value = 66157
$ws.Range("P4").Value = 6148
$ws.Range("Q4").Value = -11038
$ws.Range("R4").Value = -720
$ws.Range("S4").Value = 14341
$ws.Range("T4").Value = 145
$ws.Range("V4").Value = 118061
$ws.Range("W4").Value = 12.42
$ws.Range("X4").Value = 10.07
$ws.Range("Y4").Value = 5.25
$ws.Range("Z4").Value = 1.71
$ws.Range("AA4").Value = 231.1
$ws.Range("AB4").Value = 1024.69
$ws.Range("AC4").Value = 3016
$ws.Range("AD4").Value = 13.16
$ws.Range("AE4").Value = 60356
$ws.Range("AF4").Value = 0.66
$ws.Range("AG4").Value = 1500
$ws.Range("AH4").Value = 3.78
$ws.Range("AI4").Value = 47.05
$ws.Range("AJ4").Value = 115858891
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("U4").ClearContents()

# Row 5
$ws.Range("D5").Value = 39000
$ws.Range("E5").Value = 5056
$ws.Range("F5").Value = 5056
$ws.Range("G5").Value = 5022
$ws.Range("H5").Value = 3867
$ws.Range("I5").Value = 3867
$ws.Range("K5").Value = 230766
$ws.Range("L5").Value = 161536
$ws.Range("M5").Value = 69229
$ws.Range("N5").Value = 69229
$ws.Range("P5").Value = 6148
$ws.Range("Q5").Value = -10176
$ws.Range("R5").Value = -1309
$ws.Range("S5").Value = 11999
$ws.Range("T5").Value = 81
$ws.Range("V5").Value = 129149
$ws.Range("W5").Value = 12.96
$ws.Range("X5").Value = 9.92
$ws.Range("Y5").Value = 5.71
$ws.Range("Z5").Value = 1.72
$ws.Range("AA5").Value = 233.33
$ws.Range("AB5").Value = 1074.67
$ws.Range("AC5").Value = 3338
$ws.Range("AD5").Value = 11.86
$ws.Range("AE5").Value = 63159
$ws.Range("AF5").Value = 0.63
$ws.Range("AG5").Value = 1500
$ws.Range("AH5").Value = 3.79
$ws.Range("AI5").Value = 42.52
$ws.Range("AJ5").Value = 115858891
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()
$ws.Range("U5").ClearContents()

# Row 6
$ws.Range("D6").Value = 33542
$ws.Range("E6").Value = 4786
$ws.Range("F6").Value = 4786
$ws.Range("G6").Value = 4691
$ws.Range("H6").Value = 3452
$ws.Range("I6").Value = 3452
$ws.Range("K6").Value = 230421
$ws.Range("L6").Value = 162708
$ws.Range("M6").Value = 67713
$ws.Range("N6").Value = 67713
$ws.Range("P6").Value = 6148
$ws.Range("Q6").Value = -4031
$ws.Range("R6").Value = -705
$ws.Range("S6").Value = 6327
$ws.Range("T6").Value = 128
$ws.Range("V6").Value = 138994
$ws.Range("W6").Value = 14.27
$ws.Range("X6").Value = 10.29
$ws.Range("Y6").Value = 5.04
$ws.Range("Z6").Value = 1.5
$ws.Range("AA6").Value = 240.29
$ws.Range("AB6").Value = 1067.08
$ws.Range("AC6").Value = 2980
$ws.Range("AD6").Value = 11.59
$ws.Range("AE6").Value = 63455
$ws.Range("AF6").Value = 0.54
$ws.Range("AG6").Value = 1600
$ws.Range("AH6").Value = 4.63
$ws.Range("AI6").Value = 49.48
$ws.Range("AJ6").Value = 115858891
$ws.Range("U6").ClearContents()

# Row 7
$ws.Range("D7").Value = 33621
$ws.Range("E7").Value = 4520
$ws.Range("G7").Value = 4514
$ws.Range("H7").Value = 3455
$ws.Range("I7").Value = 3458
$ws.Range("K7").Value = 231386
$ws.Range("L7").Value = 162759
$ws.Range("M7").Value = 68630
$ws.Range("N7").Value = 53643
$ws.Range("P7").Value = 6150
$ws.Range("W7").Value = 13.44
$ws.Range("X7").Value = 10.27
$ws.Range("Y7").Value = 5.7
$ws.Range("Z7").Value = 1.5
$ws.Range("AA7").Value = 237.16
$ws.Range("AC7").Value = 2985
$ws.Range("AD7").Value = 12.85
$ws.Range("AE7").Value = 50270
$ws.Range("AF7").Value = 0.76
$ws.Range("AG7").Value = 1625
$ws.Range("AH7").Value = 4.24
$ws.Range("AI7").Value = 54.45
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()

# Row 8
$ws.Range("D8").Value = 34498
$ws.Range("E8").Value = 4968
$ws.Range("G8").Value = 4979
$ws.Range("H8").Value = 3655
$ws.Range("I8").Value = 3592
$ws.Range("K8").Value = 243510
$ws.Range("L8").Value = 172686
$ws.Range("M8").Value = 70829
$ws.Range("N8").Value = 47000
$ws.Range("P8").Value = 6150
$ws.Range("W8").Value = 14.4
$ws.Range("X8").Value = 10.59
$ws.Range("Y8").Value = 7.73
$ws.Range("Z8").Value = 1.54
$ws.Range("AA8").Value = 243.81
$ws.Range("AC8").Value = 3101
$ws.Range("AD8").Value = 12.24
$ws.Range("AE8").Value = 44044
$ws.Range("AF8").Value = 0.86
$ws.Range("AG8").Value = 1717
$ws.Range("AH8").Value = 4.52
$ws.Range("AI8").Value = 55.36
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()

# Row 9
$ws.Range("D9").Value = 35716
$ws.Range("E9").Value = 5241
$ws.Range("G9").Value = 5255
$ws.Range("H9").Value = 3865
$ws.Range("I9").Value = 3852
$ws.Range("K9").Value = 252045
$ws.Range("L9").Value = 178946
$ws.Range("M9").Value = 73098
$ws.Range("N9").Value = 48230
$ws.Range("P9").Value = 6150
$ws.Range("W9").Value = 14.67
$ws.Range("X9").Value = 10.82
$ws.Range("Y9").Value = 8.09
$ws.Range("Z9").Value = 1.56
$ws.Range("AA9").Value = 244.8
$ws.Range("AC9").Value = 3325
$ws.Range("AD9").Value = 11.41
$ws.Range("AE9").Value = 45197
$ws.Range("AF9").Value = 0.84
$ws.Range("AG9").Value = 1788
$ws.Range("AH9").Value = 4.71
$ws.Range("AI9").Value = 53.76
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
